# Append new log rows 383-398 (refactor-related log activity captured on 2025-11-28T20:22)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 383
$ws.Range("A383").Value = "2025-11-28T20:22:40.546Z"
$ws.Range("B383").Value = "HTTP"
$ws.Range("C383").Value = "GET"
$ws.Range("D383").Value = "/user/EA4C7814?roomID=109"
$ws.Range("E383").Value = 404
$ws.Range("F383").Value = "ERROR"
$ws.Range("G383").Value = ""
$ws.Range("H383").Value = ""
$ws.Range("I383").Value = "109"
$ws.Range("J383").Value = "UID EA4C7814 not found in /user"
$ws.Range("K383").Value = "::ffff:172.28.219.204"

# Row 384
$ws.Range("A384").Value = "2025-11-28T20:22:40.877Z"
$ws.Range("B384").Value = "HTTP"
$ws.Range("C384").Value = "GET"
$ws.Range("D384").Value = "/uid-name/EA4C7814"
$ws.Range("E384").Value = 404
$ws.Range("F384").Value = "ERROR"
$ws.Range("G384").Value = ""
$ws.Range("H384").Value = ""
$ws.Range("I384").Value = ""
$ws.Range("J384").Value = "UID EA4C7814 not found in /uid-name"
$ws.Range("K384").Value = "::ffff:172.28.219.204"

# Row 385
$ws.Range("A385").Value = "2025-11-28T20:22:42.653Z"
$ws.Range("B385").Value = "HTTP"
$ws.Range("C385").Value = "GET"
$ws.Range("D385").Value = "/user/EA4C7814?roomID=110"
$ws.Range("E385").Value = 404
$ws.Range("F385").Value = "ERROR"
$ws.Range("G385").Value = ""
$ws.Range("H385").Value = ""
$ws.Range("I385").Value = "110"
$ws.Range("J385").Value = "UID EA4C7814 not found in /user"
$ws.Range("K385").Value = "::ffff:172.28.219.204"

# Row 386
$ws.Range("A386").Value = "2025-11-28T20:22:42.729Z"
$ws.Range("B386").Value = "HTTP"
$ws.Range("C386").Value = "GET"
$ws.Range("D386").Value = "/uid-name/EA4C7814"
$ws.Range("E386").Value = 404
$ws.Range("F386").Value = "ERROR"
$ws.Range("G386").Value = ""
$ws.Range("H386").Value = ""
$ws.Range("I386").Value = ""
$ws.Range("J386").Value = "UID EA4C7814 not found in /uid-name"
$ws.Range("K386").Value = "::ffff:172.28.219.204"

# Row 387
$ws.Range("A387").Value = "2025-11-28T20:22:47.872Z"
$ws.Range("B387").Value = "HTTP"
$ws.Range("C387").Value = "GET"
$ws.Range("D387").Value = "/user/EA4C7814?roomID=107"
$ws.Range("E387").Value = 404
$ws.Range("F387").Value = "ERROR"
$ws.Range("G387").Value = ""
$ws.Range("H387").Value = ""
$ws.Range("I387").Value = "107"
$ws.Range("J387").Value = "UID EA4C7814 not found in /user"
$ws.Range("K387").Value = "::ffff:172.28.219.204"

# Row 388
$ws.Range("A388").Value = "2025-11-28T20:22:48.183Z"
$ws.Range("B388").Value = "HTTP"
$ws.Range("C388").Value = "GET"
$ws.Range("D388").Value = "/uid-name/EA4C7814"
$ws.Range("E388").Value = 404
$ws.Range("F388").Value = "ERROR"
$ws.Range("G388").Value = ""
$ws.Range("H388").Value = ""
$ws.Range("I388").Value = ""
$ws.Range("J388").Value = "UID EA4C7814 not found in /uid-name"
$ws.Range("K388").Value = "::ffff:172.28.219.204"

# Row 389
$ws.Range("A389").Value = "2025-11-28T20:22:48.696Z"
$ws.Range("B389").Value = "HTTP"
$ws.Range("C389").Value = "GET"
$ws.Range("D389").Value = "/user/EA4C7814?roomID=107"
$ws.Range("E389").Value = 404
$ws.Range("F389").Value = "ERROR"
$ws.Range("G389").Value = ""
$ws.Range("H389").Value = ""
$ws.Range("I389").Value = "107"
$ws.Range("J389").Value = "UID EA4C7814 not found in /user"
$ws.Range("K389").Value = "::ffff:172.28.219.204"

# Row 390
$ws.Range("A390").Value = "2025-11-28T20:22:49.137Z"
$ws.Range("B390").Value = "HTTP"
$ws.Range("C390").Value = "GET"
$ws.Range("D390").Value = "/uid-name/EA4C7814"
$ws.Range("E390").Value = 404
$ws.Range("F390").Value = "ERROR"
$ws.Range("G390").Value = ""
$ws.Range("H390").Value = ""
$ws.Range("I390").Value = ""
$ws.Range("J390").Value = "UID EA4C7814 not found in /uid-name"
$ws.Range("K390").Value = "::ffff:172.28.219.204"

# Row 391
$ws.Range("A391").Value = "2025-11-28T20:22:51.265Z"
$ws.Range("B391").Value = "HTTP"
$ws.Range("C391").Value = "GET"
$ws.Range("D391").Value = "/user/EA4C7814?roomID=106"
$ws.Range("E391").Value = 404
$ws.Range("F391").Value = "ERROR"
$ws.Range("G391").Value = ""
$ws.Range("H391").Value = ""
$ws.Range("I391").Value = "106"
$ws.Range("J391").Value = "UID EA4C7814 not found in /user"
$ws.Range("K391").Value = "::ffff:172.28.219.204"

# Row 392
$ws.Range("A392").Value = "2025-11-28T20:22:51.338Z"
$ws.Range("B392").Value = "HTTP"
$ws.Range("C392").Value = "GET"
$ws.Range("D392").Value = "/uid-name/EA4C7814"
$ws.Range("E392").Value = 404
$ws.Range("F392").Value = "ERROR"
$ws.Range("G392").Value = ""
$ws.Range("H392").Value = ""
$ws.Range("I392").Value = ""
$ws.Range("J392").Value = "UID EA4C7814 not found in /uid-name"
$ws.Range("K392").Value = "::ffff:172.28.219.204"

# Row 393
$ws.Range("A393").Value = "2025-11-28T20:22:53.820Z"
$ws.Range("B393").Value = "HTTP"
$ws.Range("C393").Value = "GET"
$ws.Range("D393").Value = "/user/EA4C7814?roomID=105"
$ws.Range("E393").Value = 404
$ws.Range("F393").Value = "ERROR"
$ws.Range("G393").Value = ""
$ws.Range("H393").Value = ""
$ws.Range("I393").Value = "105"
$ws.Range("J393").Value = "UID EA4C7814 not found in /user"
$ws.Range("K393").Value = "::ffff:172.28.219.204"

# Row 394
$ws.Range("A394").Value = "2025-11-28T20:22:53.951Z"
$ws.Range("B394").Value = "HTTP"
$ws.Range("C394").Value = "GET"
$ws.Range("D394").Value = "/uid-name/EA4C7814"
$ws.Range("E394").Value = 404
$ws.Range("F394").Value = "ERROR"
$ws.Range("G394").Value = ""
$ws.Range("H394").Value = ""
$ws.Range("I394").Value = ""
$ws.Range("J394").Value = "UID EA4C7814 not found in /uid-name"
$ws.Range("K394").Value = "::ffff:172.28.219.204"

# Row 395
$ws.Range("A395").Value = "2025-11-28T20:22:54.042Z"
$ws.Range("B395").Value = "HTTP"
$ws.Range("C395").Value = "GET"
$ws.Range("D395").Value = "/user/EA4C7814?roomID=106"
$ws.Range("E395").Value = 404
$ws.Range("F395").Value = "ERROR"
$ws.Range("G395").Value = ""
$ws.Range("H395").Value = ""
$ws.Range("I395").Value = "106"
$ws.Range("J395").Value = "UID EA4C7814 not found in /user"
$ws.Range("K395").Value = "::ffff:172.28.219.204"

# Row 396
$ws.Range("A396").Value = "2025-11-28T20:22:54.092Z"
$ws.Range("B396").Value = "HTTP"
$ws.Range("C396").Value = "GET"
$ws.Range("D396").Value = "/uid-name/EA4C7814"
$ws.Range("E396").Value = 404
$ws.Range("F396").Value = "ERROR"
$ws.Range("G396").Value = ""
$ws.Range("H396").Value = ""
$ws.Range("I396").Value = ""
$ws.Range("J396").Value = "UID EA4C7814 not found in /uid-name"
$ws.Range("K396").Value = "::ffff:172.28.219.204"

# Row 397
$ws.Range("A397").Value = "2025-11-28T20:22:55.777Z"
$ws.Range("B397").Value = "HTTP"
$ws.Range("C397").Value = "GET"
$ws.Range("D397").Value = "/user/EA4C7814?roomID=105"
$ws.Range("E397").Value = 404
$ws.Range("F397").Value = "ERROR"
$ws.Range("G397").Value = ""
$ws.Range("H397").Value = ""
$ws.Range("I397").Value = "105"
$ws.Range("J397").Value = "UID EA4C7814 not found in /user"
$ws.Range("K397").Value = "::ffff:172.28.219.204"

# Row 398
$ws.Range("A398").Value = "2025-11-28T20:22:55.829Z"
$ws.Range("B398").Value = "HTTP"
$ws.Range("C398").Value = "GET"
$ws.Range("D398").Value = "/uid-name/EA4C7814"
$ws.Range("E398").Value = 404
$ws.Range("F398").Value = "ERROR"
$ws.Range("G398").Value = ""
$ws.Range("H398").Value = ""
$ws.Range("I398").Value = ""
$ws.Range("J398").Value = "UID EA4C7814 not found in /uid-name"
$ws.Range("K398").Value = "::ffff:172.28.219.204"

